# Apply crypto price/volume updates per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.832.51"
$ws.Range("E2").Value = "  +0.08%  "
$ws.Range("D3").Value = "3.787.51"
$ws.Range("E3").Value = "  -0.64%  "
$ws.Range("D4").Value = "'0.995"
$ws.Range("E4").Value = "  -0.52%  "
$ws.Range("D5").Value = "'602.59"
$ws.Range("E5").Value = "  -0.88%  "
$ws.Range("D6").Value = "'163.12"
$ws.Range("E6").Value = "  -2.32%  "
$ws.Range("D7").Value = "3.785.89"
$ws.Range("E7").Value = "  -0.63%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("D10").Value = "'0.156"
$ws.Range("E10").Value = "  -2.86%  "
$ws.Range("E11").Value = "  -1.09%  "
$ws.Range("D12").Value = "'6.73"
$ws.Range("E12").Value = "  +6.80%  "
$ws.Range("E13").Value = "  -3.56%  "
$ws.Range("D14").Value = "'35.14"
$ws.Range("E14").Value = "  -2.57%  "
$ws.Range("D15").Value = "4.425.38"
$ws.Range("E15").Value = "  -0.51%  "
$ws.Range("D16").Value = "3.802.67"
$ws.Range("E16").Value = "  +0.02%  "
$ws.Range("D17").Value = "67.825.13"
$ws.Range("E17").Value = "  +0.04%  "
$ws.Range("E18").Value = "  -1.95%  "
$ws.Range("E19").Value = "  +2.01%  "
$ws.Range("E20").Value = "  -1.34%  "
$ws.Range("D21").Value = "'457.59"
$ws.Range("E21").Value = "  -1.07%  "
$ws.Range("D22").Value = "'9.45"
$ws.Range("E22").Value = "  -4.70%  "
$ws.Range("E23").Value = "  -1.53%  "
$ws.Range("D24").Value = "'82.96"
$ws.Range("E24").Value = "  -0.59%  "
$ws.Range("D25").Value = "'0.0000142"
$ws.Range("E25").Value = "  -5.94%  "
$ws.Range("D26").Value = "'11.84"
$ws.Range("E26").Value = "  -2.09%  "
$ws.Range("D27").Value = "'2.07"
$ws.Range("E27").Value = "  -1.80%  "
$ws.Range("D29").Value = "'9.90"
$ws.Range("E29").Value = "  -1.49%  "
$ws.Range("D30").Value = "3.939.04"
$ws.Range("E30").Value = "  -0.51%  "
$ws.Range("D31").Value = "'7.21"
$ws.Range("E31").Value = "  -2.87%  "
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").Value = "'2.59"
$ws.Range("E32").Value = "  -7.56%  "
$ws.Range("B33").Value = "ImmutableX"
$ws.Range("C33").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D33").Value = "'2.18"
$ws.Range("E33").Value = "  -2.32%  "
$ws.Range("D34").Value = "'28.91"
$ws.Range("E34").Value = "  -2.44%  "
$ws.Range("E35").Value = "  +0.56%  "
$ws.Range("E36").Value = "  -1.84%  "
$ws.Range("D37").Value = "'0.0990"
$ws.Range("E37").Value = "  -1.14%  "
$ws.Range("D38").Value = "'0.144"
$ws.Range("E38").Value = "  +4.76%  "
$ws.Range("D39").Value = "'5.79"
$ws.Range("E39").Value = "  -0.32%  "
$ws.Range("D40").Value = "'0.978"
$ws.Range("E40").Value = "  -2.27%  "
$ws.Range("D41").Value = "'3.17"
$ws.Range("E41").Value = "  -6.62%  "
$ws.Range("D42").Value = "'1.00"
$ws.Range("E42").Value = "  +0.05%  "
$ws.Range("D44").Value = "'43.65"
$ws.Range("D45").Value = "'47.17"
$ws.Range("E45").Value = "  -2.05%  "
$ws.Range("D46").Value = "'151.89"
$ws.Range("E46").Value = "  +2.06%  "
$ws.Range("E47").Value = "  -2.47%  "
$ws.Range("D48").Value = "'8.27"
$ws.Range("E48").Value = "  -1.09%  "
$ws.Range("E49").Value = "  -0.86%  "
$ws.Range("E50").Value = "  -1.08%  "
$ws.Range("D51").Value = "'26.47"
$ws.Range("E51").Value = "  -6.11%  "
